$d = $word.ActiveDocument
$failures = @()

if (-not $d.Content.Find.Execute("2023-02-02 Thursday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-02-03 Friday", 2)) { $failures += "2023-02-02 Thursday" }
if (-not $d.Content.Find.Execute("82-36=", $true, $false, $false, $false, $false, $true, 1, $false, "60+10=", 2)) { $failures += "82-36=" }
if (-not $d.Content.Find.Execute("82-54=", $true, $false, $false, $false, $false, $true, 1, $false, "94-4=", 2)) { $failures += "82-54=" }
if (-not $d.Content.Find.Execute("51+9=", $true, $false, $false, $false, $false, $true, 1, $false, "14+65=", 2)) { $failures += "51+9=" }
if (-not $d.Content.Find.Execute("65-17=", $true, $false, $false, $false, $false, $true, 1, $false, "86-29=", 2)) { $failures += "65-17=" }
if (-not $d.Content.Find.Execute("54-9=", $true, $false, $false, $false, $false, $true, 1, $false, "53+46=", 2)) { $failures += "54-9=" }
if (-not $d.Content.Find.Execute("56+1=", $true, $false, $false, $false, $false, $true, 1, $false, "25+27=", 2)) { $failures += "56+1=" }
if (-not $d.Content.Find.Execute("82+13=", $true, $false, $false, $false, $false, $true, 1, $false, "82-27=", 2)) { $failures += "82+13=" }
if (-not $d.Content.Find.Execute("75-50=", $true, $false, $false, $false, $false, $true, 1, $false, "80-23=", 2)) { $failures += "75-50=" }
if (-not $d.Content.Find.Execute("27-3=", $true, $false, $false, $false, $false, $true, 1, $false, "89-71=", 2)) { $failures += "27-3=" }
if (-not $d.Content.Find.Execute("97-34=", $true, $false, $false, $false, $false, $true, 1, $false, "58-52=", 2)) { $failures += "97-34=" }
if (-not $d.Content.Find.Execute("95-38=", $true, $false, $false, $false, $false, $true, 1, $false, "0+32=", 2)) { $failures += "95-38=" }
if (-not $d.Content.Find.Execute("92-46=", $true, $false, $false, $false, $false, $true, 1, $false, "65+9=", 2)) { $failures += "92-46=" }
if (-not $d.Content.Find.Execute("39+54=", $true, $false, $false, $false, $false, $true, 1, $false, "23+31=", 2)) { $failures += "39+54=" }
if (-not $d.Content.Find.Execute("56-15=", $true, $false, $false, $false, $false, $true, 1, $false, "43+50=", 2)) { $failures += "56-15=" }
if (-not $d.Content.Find.Execute("5-3=", $true, $false, $false, $false, $false, $true, 1, $false, "20+14=", 2)) { $failures += "5-3=" }
if (-not $d.Content.Find.Execute("3+24=", $true, $false, $false, $false, $false, $true, 1, $false, "14+34=", 2)) { $failures += "3+24=" }
if (-not $d.Content.Find.Execute("47+8=", $true, $false, $false, $false, $false, $true, 1, $false, "68+23=", 2)) { $failures += "47+8=" }
if (-not $d.Content.Find.Execute("31+56=", $true, $false, $false, $false, $false, $true, 1, $false, "61+0=", 2)) { $failures += "31+56=" }
if (-not $d.Content.Find.Execute("71-6=", $true, $false, $false, $false, $false, $true, 1, $false, "83-60=", 2)) { $failures += "71-6=" }
if (-not $d.Content.Find.Execute("3+9=", $true, $false, $false, $false, $false, $true, 1, $false, "36+33=", 2)) { $failures += "3+9=" }
if (-not $d.Content.Find.Execute("57-55=", $true, $false, $false, $false, $false, $true, 1, $false, "26+70=", 2)) { $failures += "57-55=" }
if (-not $d.Content.Find.Execute("89-65=", $true, $false, $false, $false, $false, $true, 1, $false, "71-69=", 2)) { $failures += "89-65=" }
if (-not $d.Content.Find.Execute("55+9=", $true, $false, $false, $false, $false, $true, 1, $false, "55-19=", 2)) { $failures += "55+9=" }
if (-not $d.Content.Find.Execute("50+23=", $true, $false, $false, $false, $false, $true, 1, $false, "62-47=", 2)) { $failures += "50+23=" }
if (-not $d.Content.Find.Execute("90+1=", $true, $false, $false, $false, $false, $true, 1, $false, "90-64=", 2)) { $failures += "90+1=" }
if (-not $d.Content.Find.Execute("58-26=", $true, $false, $false, $false, $false, $true, 1, $false, "27+41=", 2)) { $failures += "58-26=" }
if (-not $d.Content.Find.Execute("5+85=", $true, $false, $false, $false, $false, $true, 1, $false, "73-48=", 2)) { $failures += "5+85=" }
if (-not $d.Content.Find.Execute("68-4=", $true, $false, $false, $false, $false, $true, 1, $false, "90-34=", 2)) { $failures += "68-4=" }
if (-not $d.Content.Find.Execute("83-31=", $true, $false, $false, $false, $false, $true, 1, $false, "43+42=", 2)) { $failures += "83-31=" }
if (-not $d.Content.Find.Execute("93-74=", $true, $false, $false, $false, $false, $true, 1, $false, "41-23=", 2)) { $failures += "93-74=" }
if (-not $d.Content.Find.Execute("94-72=", $true, $false, $false, $false, $false, $true, 1, $false, "29+35=", 2)) { $failures += "94-72=" }
if (-not $d.Content.Find.Execute("44+19=", $true, $false, $false, $false, $false, $true, 1, $false, "19+61=", 2)) { $failures += "44+19=" }
if (-not $d.Content.Find.Execute("61+3=", $true, $false, $false, $false, $false, $true, 1, $false, "87-79=", 2)) { $failures += "61+3=" }
if (-not $d.Content.Find.Execute("27+42=", $true, $false, $false, $false, $false, $true, 1, $false, "63-13=", 2)) { $failures += "27+42=" }
if (-not $d.Content.Find.Execute("52-3=", $true, $false, $false, $false, $false, $true, 1, $false, "4+30=", 2)) { $failures += "52-3=" }
if (-not $d.Content.Find.Execute("20+32=", $true, $false, $false, $false, $false, $true, 1, $false, "54-13=", 2)) { $failures += "20+32=" }
if (-not $d.Content.Find.Execute("54-49=", $true, $false, $false, $false, $false, $true, 1, $false, "95-24=", 2)) { $failures += "54-49=" }
if (-not $d.Content.Find.Execute("40+56=", $true, $false, $false, $false, $false, $true, 1, $false, "58+32=", 2)) { $failures += "40+56=" }
if (-not $d.Content.Find.Execute("38-9=", $true, $false, $false, $false, $false, $true, 1, $false, "6+70=", 2)) { $failures += "38-9=" }
if (-not $d.Content.Find.Execute("70-67=", $true, $false, $false, $false, $false, $true, 1, $false, "68-14=", 2)) { $failures += "70-67=" }
if (-not $d.Content.Find.Execute("67-39=", $true, $false, $false, $false, $false, $true, 1, $false, "41-1=", 2)) { $failures += "67-39=" }
if (-not $d.Content.Find.Execute("89-74=", $true, $false, $false, $false, $false, $true, 1, $false, "54+14=", 2)) { $failures += "89-74=" }
if (-not $d.Content.Find.Execute("26+51=", $true, $false, $false, $false, $false, $true, 1, $false, "0+29=", 2)) { $failures += "26+51=" }
if (-not $d.Content.Find.Execute("56+42=", $true, $false, $false, $false, $false, $true, 1, $false, "29+9=", 2)) { $failures += "56+42=" }
if (-not $d.Content.Find.Execute("86-81=", $true, $false, $false, $false, $false, $true, 1, $false, "49-16=", 2)) { $failures += "86-81=" }
if (-not $d.Content.Find.Execute("50+34=", $true, $false, $false, $false, $false, $true, 1, $false, "39+2=", 2)) { $failures += "50+34=" }
if (-not $d.Content.Find.Execute("27+59=", $true, $false, $false, $false, $false, $true, 1, $false, "44+54=", 2)) { $failures += "27+59=" }
if (-not $d.Content.Find.Execute("55+3=", $true, $false, $false, $false, $false, $true, 1, $false, "27+1=", 2)) { $failures += "55+3=" }
if (-not $d.Content.Find.Execute("63+8=", $true, $false, $false, $false, $false, $true, 1, $false, "85-55=", 2)) { $failures += "63+8=" }
if (-not $d.Content.Find.Execute("5+3=", $true, $false, $false, $false, $false, $true, 1, $false, "86+8=", 2)) { $failures += "5+3=" }
if (-not $d.Content.Find.Execute("1+17=", $true, $false, $false, $false, $false, $true, 1, $false, "80-49=", 2)) { $failures += "1+17=" }
if (-not $d.Content.Find.Execute("67-41=", $true, $false, $false, $false, $false, $true, 1, $false, "3+28=", 2)) { $failures += "67-41=" }
if (-not $d.Content.Find.Execute("49-29=", $true, $false, $false, $false, $false, $true, 1, $false, "17-3=", 2)) { $failures += "49-29=" }
if (-not $d.Content.Find.Execute("25+66=", $true, $false, $false, $false, $false, $true, 1, $false, "42+46=", 2)) { $failures += "25+66=" }
if (-not $d.Content.Find.Execute("30+35=", $true, $false, $false, $false, $false, $true, 1, $false, "57-1=", 2)) { $failures += "30+35=" }
if (-not $d.Content.Find.Execute("66+14=", $true, $false, $false, $false, $false, $true, 1, $false, "64+26=", 2)) { $failures += "66+14=" }
if (-not $d.Content.Find.Execute("84-39=", $true, $false, $false, $false, $false, $true, 1, $false, "13+4=", 2)) { $failures += "84-39=" }
if (-not $d.Content.Find.Execute("52+6=", $true, $false, $false, $false, $false, $true, 1, $false, "88+11=", 2)) { $failures += "52+6=" }
if (-not $d.Content.Find.Execute("9+47=", $true, $false, $false, $false, $false, $true, 1, $false, "55+23=", 2)) { $failures += "9+47=" }
if (-not $d.Content.Find.Execute("43+21=", $true, $false, $false, $false, $false, $true, 1, $false, "86-22=", 2)) { $failures += "43+21=" }
if (-not $d.Content.Find.Execute("32-3=", $true, $false, $false, $false, $false, $true, 1, $false, "74-34=", 2)) { $failures += "32-3=" }
if (-not $d.Content.Find.Execute("53+30=", $true, $false, $false, $false, $false, $true, 1, $false, "55+34=", 2)) { $failures += "53+30=" }
if (-not $d.Content.Find.Execute("37+23=", $true, $false, $false, $false, $false, $true, 1, $false, "22+24=", 2)) { $failures += "37+23=" }
if (-not $d.Content.Find.Execute("15+0=", $true, $false, $false, $false, $false, $true, 1, $false, "71-68=", 2)) { $failures += "15+0=" }
if (-not $d.Content.Find.Execute("75-12=", $true, $false, $false, $false, $false, $true, 1, $false, "51-2=", 2)) { $failures += "75-12=" }
if (-not $d.Content.Find.Execute("68+0=", $true, $false, $false, $false, $false, $true, 1, $false, "23+30=", 2)) { $failures += "68+0=" }
if (-not $d.Content.Find.Execute("62-20=", $true, $false, $false, $false, $false, $true, 1, $false, "9+58=", 2)) { $failures += "62-20=" }
if (-not $d.Content.Find.Execute("53-34=", $true, $false, $false, $false, $false, $true, 1, $false, "63-33=", 2)) { $failures += "53-34=" }
if (-not $d.Content.Find.Execute("56+38=", $true, $false, $false, $false, $false, $true, 1, $false, "3+31=", 2)) { $failures += "56+38=" }
if (-not $d.Content.Find.Execute("45+52=", $true, $false, $false, $false, $false, $true, 1, $false, "24+54=", 2)) { $failures += "45+52=" }
if (-not $d.Content.Find.Execute("98-22=", $true, $false, $false, $false, $false, $true, 1, $false, "17+24=", 2)) { $failures += "98-22=" }
if (-not $d.Content.Find.Execute("5+13=", $true, $false, $false, $false, $false, $true, 1, $false, "30+38=", 2)) { $failures += "5+13=" }
if (-not $d.Content.Find.Execute("41-11=", $true, $false, $false, $false, $false, $true, 1, $false, "15+13=", 2)) { $failures += "41-11=" }
if (-not $d.Content.Find.Execute("78-32=", $true, $false, $false, $false, $false, $true, 1, $false, "84+0=", 2)) { $failures += "78-32=" }
if (-not $d.Content.Find.Execute("26+73=", $true, $false, $false, $false, $false, $true, 1, $false, "64-1=", 2)) { $failures += "26+73=" }
if (-not $d.Content.Find.Execute("67+9=", $true, $false, $false, $false, $false, $true, 1, $false, "43-28=", 2)) { $failures += "67+9=" }
if (-not $d.Content.Find.Execute("7+22=", $true, $false, $false, $false, $false, $true, 1, $false, "75-69=", 2)) { $failures += "7+22=" }
if (-not $d.Content.Find.Execute("36-0=", $true, $false, $false, $false, $false, $true, 1, $false, "17-2=", 2)) { $failures += "36-0=" }
if (-not $d.Content.Find.Execute("13+45=", $true, $false, $false, $false, $false, $true, 1, $false, "90-12=", 2)) { $failures += "13+45=" }
if (-not $d.Content.Find.Execute("25+59=", $true, $false, $false, $false, $false, $true, 1, $false, "54+22=", 2)) { $failures += "25+59=" }
if (-not $d.Content.Find.Execute("9+78=", $true, $false, $false, $false, $false, $true, 1, $false, "29+45=", 2)) { $failures += "9+78=" }
if (-not $d.Content.Find.Execute("19+56=", $true, $false, $false, $false, $false, $true, 1, $false, "48+37=", 2)) { $failures += "19+56=" }
if (-not $d.Content.Find.Execute("16+5=", $true, $false, $false, $false, $false, $true, 1, $false, "0+75=", 2)) { $failures += "16+5=" }
if (-not $d.Content.Find.Execute("88-64=", $true, $false, $false, $false, $false, $true, 1, $false, "72+19=", 2)) { $failures += "88-64=" }
if (-not $d.Content.Find.Execute("32-26=", $true, $false, $false, $false, $false, $true, 1, $false, "56-34=", 2)) { $failures += "32-26=" }
if (-not $d.Content.Find.Execute("41-12=", $true, $false, $false, $false, $false, $true, 1, $false, "98-30=", 2)) { $failures += "41-12=" }
if (-not $d.Content.Find.Execute("61-41=", $true, $false, $false, $false, $false, $true, 1, $false, "3+95=", 2)) { $failures += "61-41=" }
if (-not $d.Content.Find.Execute("26-0=", $true, $false, $false, $false, $false, $true, 1, $false, "81-27=", 2)) { $failures += "26-0=" }
if (-not $d.Content.Find.Execute("44-7=", $true, $false, $false, $false, $false, $true, 1, $false, "25-6=", 2)) { $failures += "44-7=" }
if (-not $d.Content.Find.Execute("18+38=", $true, $false, $false, $false, $false, $true, 1, $false, "30-2=", 2)) { $failures += "18+38=" }
if (-not $d.Content.Find.Execute("39+29=", $true, $false, $false, $false, $false, $true, 1, $false, "90-63=", 2)) { $failures += "39+29=" }
if (-not $d.Content.Find.Execute("74+22=", $true, $false, $false, $false, $false, $true, 1, $false, "39+29=", 2)) { $failures += "74+22=" }
if (-not $d.Content.Find.Execute("91-72=", $true, $false, $false, $false, $false, $true, 1, $false, "55-37=", 2)) { $failures += "91-72=" }
if (-not $d.Content.Find.Execute("11+28=", $true, $false, $false, $false, $false, $true, 1, $false, "91-27=", 2)) { $failures += "11+28=" }
if (-not $d.Content.Find.Execute("43-34=", $true, $false, $false, $false, $false, $true, 1, $false, "3+6=", 2)) { $failures += "43-34=" }
if (-not $d.Content.Find.Execute("63-30=", $true, $false, $false, $false, $false, $true, 1, $false, "61-46=", 2)) { $failures += "63-30=" }
if (-not $d.Content.Find.Execute("95-65=", $true, $false, $false, $false, $false, $true, 1, $false, "13+75=", 2)) { $failures += "95-65=" }
if (-not $d.Content.Find.Execute("18+22=", $true, $false, $false, $false, $false, $true, 1, $false, "81-19=", 2)) { $failures += "18+22=" }
if (-not $d.Content.Find.Execute("29-0=", $true, $false, $false, $false, $false, $true, 1, $false, "65-56=", 2)) { $failures += "29-0=" }
if (-not $d.Content.Find.Execute("70-65=", $true, $false, $false, $false, $false, $true, 1, $false, "2+41=", 2)) { $failures += "70-65=" }

if ($failures.Count -gt 0) { "FAILED: " + ($failures -join ", ") } else { "OK: all replacements applied" }
